# This edit corresponds to the commit "alteração de indexador para recursos neutros":
# Rows whose "NomePrecoTransferencia" (column O) pointed to "IPCA Mensal" now point
# to "DI Mensal" instead, and their "CodPrecoTransferencia" (column N) changes from 1 to 4
# (the same code already used by the other rows referencing "DI Mensal").
# Because no cell references "IPCA Mensal" any more after this change, the shared string
# is dropped from the workbook's string table on save, which is what shifts every other
# index (>= 101) down by one in the underlying XML - that happens automatically, we only
# need to change the cell contents below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows (1-indexed, matching worksheet row numbers) where CodPrecoTransferencia (col N)
# was 1 and NomePrecoTransferencia (col O) was "IPCA Mensal".
$rows = @(
    4,5,6,7,8,9,10,11,12,
    132,133,134,135,136,137,138,139,
    145,146,147,148,149,150,151,152,153,154,155,156,157,158,159,160,161,162,163,164,165,166,167,
    525,526,527,528,529,530,531,532,533,534,535,536,537,538,539,540,541,542,543,544,545,546,
    548
)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 14).Value2 = 4
    $ws.Cells.Item($r, 15).Value2 = "DI Mensal"
}
